$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" date column (C) for rows 2-5 from 2023-10-09 to 2023-10-13
# (Excel serial date 45208 -> 45212), keeping existing cell formatting.
$newDate = [DateTime]::FromOADate(45212)

$ws.Range("C2").Value = $newDate
$ws.Range("C3").Value = $newDate
$ws.Range("C4").Value = $newDate
$ws.Range("C5").Value = $newDate
